# Excel template bug fix
# Update the "Obs_relatorio" (column F) messages on the "Tanque" sheet for
# rows 2-4 so they correctly report the divergence between the SPED value
# and the report value instead of always claiming a successful validation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tanque")

$ws.Range("F2").Value = "Divergência entre o SPED(954,00) e o relatório(40790,95)!"
$ws.Range("F3").Value = "Divergência entre o SPED(954,00) e o relatório(206100,72)!"
$ws.Range("F4").Value = "Divergência entre o SPED(954,00) e o relatório(68167,68)!"
